# Re-create the "delanalyse1_" sheet: a second sheet that mirrors the
# used range (A1:G5) of the original "delanalyse1_noegletalniveau1og5"
# sheet - same values, number formats and borders - and becomes the
# active/selected sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "delanalyse1_"

# Copy values then formats (two passes) so both land on the new sheet.
$ws1.Range("A1:G5").Copy()
$ws2.Range("A1:G5").PasteSpecial(-4163)   # xlPasteValues
$ws1.Range("A1:G5").Copy()
$ws2.Range("A1:G5").PasteSpecial(-4122)   # xlPasteFormats

# Slightly widen column D on the original sheet.
$ws1.Columns.Item(4).ColumnWidth = 12.43

# Selections: sheet1 no longer the active/selected tab, new cell picked;
# sheet2 (now active) gets its own selection.
$ws1.Range("G32").Select() | Out-Null
$ws2.Range("G12").Select() | Out-Null

Write-Host "delanalyse1_ sheet added"
